# "update 5 jul bis" — append the 2021-06-30 (serial 44377) indicator rows
# for Andalucía and its eight provinces to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the date style already used in column A (numFmtId 165, "yyyy-mm-dd").
$ws.Range("A263:A271").NumberFormat = "yyyy\-mm\-dd"

$rows = @(
    @{ Row = 263; Territorio = "Andalucía"; C = 166.4; D = 78.4;  E = 38;   F = 17.8; G = 8.1;  H = 25.3; I = 2.4; J = 7;    Riesgo = "BAJO"  },
    @{ Row = 264; Territorio = "Almería";   C = 76.8;  D = 41.3;  E = 20.5; F = 16;   G = 5.5;  H = 44.5; I = 0.9; J = 2.3;  Riesgo = "BAJO"  },
    @{ Row = 265; Territorio = "Cádiz";     C = 145.3; D = 69.7;  E = 38.4; F = 16.3; G = 8.3;  H = 30.2; I = 1.7; J = 3.3;  Riesgo = "BAJO"  },
    @{ Row = 266; Territorio = "Córdoba";   C = 247.5; D = 100.6; E = 76.6; F = 37;   G = 10.3; H = 27.9; I = 3.3; J = 12.7; Riesgo = "MEDIO" },
    @{ Row = 267; Territorio = "Granada";   C = 217.7; D = 106.3; E = 51;   F = 26.4; G = 11.4; H = 8.6;  I = 2.7; J = 7.4;  Riesgo = "BAJO"  },
    @{ Row = 268; Territorio = "Huelva";    C = 161.2; D = 70;    E = 29.5; F = 13.6; G = 6.3;  H = 53.7; I = 2.4; J = 5.5;  Riesgo = "BAJO"  },
    @{ Row = 269; Territorio = "Jaén";      C = 150.1; D = 79.3;  E = 28.5; F = 4.1;  G = 9.5;  H = 35.1; I = 4.3; J = 9.8;  Riesgo = "BAJO"  },
    @{ Row = 270; Territorio = "Málaga";    C = 180.4; D = 93.5;  E = 31.4; F = 17.1; G = 10.1; H = 21.1; I = 2.3; J = 4;    Riesgo = "BAJO"  },
    @{ Row = 271; Territorio = "Sevilla";   C = 149.7; D = 64.4;  E = 31.4; F = 12.9; G = 8.1;  H = 21.3; I = 2.5; J = 11.7; Riesgo = "MEDIO" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 44377
    $ws.Cells.Item($row, 2).Value = $r.Territorio
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.Riesgo
}

# Move selection to the new "next empty" data-entry cell, mirroring the
# author's saved cursor position after appending the block.
$ws.Range("C272").Select()
